$wb = $excel.ActiveWorkbook

# --- 1) Update status text "Ready for handoff" -> "In Translation" everywhere it occurs ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $cellText = [string]$cell.Value2
            if ($cellText -eq "Ready for handoff") {
                $cell.Value = "In Translation"
            }
        }
    }
}

# --- 2) Narrow the Status columns to match the shorter text ---
# "Overview" sheet: status columns are E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.576851254417766
$overview.Columns.Item(6).ColumnWidth = 12.576851254417766

# "zh-cn" and "de-de" sheets: status column is C
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.576851254417766

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.576851254417766
